$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.611.23"
$ws.Range("E2").Value = "  -2.99%  "

$ws.Range("D3").Value = "3.809.80"
$ws.Range("E3").Value = "  +1.84%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.93%  "

$ws.Range("D7").Value = "3.807.94"
$ws.Range("E7").Value = "  +1.82%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.522"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.32%  "

$ws.Range("E10").Value = "  -4.10%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.20"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.42%  "

$ws.Range("E14").Value = "  -3.87%  "

$ws.Range("D15").Value = "4.433.70"
$ws.Range("E15").Value = "  +1.62%  "

$ws.Range("D16").Value = "3.792.70"
$ws.Range("E16").Value = "  +1.44%  "

$ws.Range("D17").Value = "67.700.79"
$ws.Range("E17").Value = "  -2.98%  "

$ws.Range("E18").Value = "  -4.72%  "

$ws.Range("E19").Value = "  -4.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.83%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "487.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.96%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.724"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.44%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -11.62%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000140"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.36%  "

$ws.Range("E27").Value = "  -5.70%  "

$ws.Range("E28").Value = "  -13.75%  "

$ws.Range("E29").Value = "  -0.15%  "

$ws.Range("E30").Value = "  -0.31%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.67"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.57%  "

$ws.Range("E32").Value = "  -3.18%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.80"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.31%  "

$ws.Range("E34").Value = "  -4.86%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.16%  "

$ws.Range("E36").Value = "  -3.99%  "

$ws.Range("E38").Value = "  -6.32%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "451.08"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.44%  "

$ws.Range("E40").Value = "  -9.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "48.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.65%  "

$ws.Range("E42").Value = "  -3.96%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.51%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.51%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -11.34%  "

$ws.Range("D46").Value = "2.838.67"
$ws.Range("E46").Value = "  -3.95%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.62%  "

$ws.Range("E48").Value = "  +0.01%  "

$ws.Range("E49").Value = "  -3.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.36%  "
